# "demo pva with notes"
#
# The slide-12 callout box ("Rectangle 4") currently reads:
#   "In 20% of the replications, all Coastal  populations ended in the
#    Extirpated state"
# and needs to become:
#   "In 21% of the replications, all Coastal  populations ended in the
#    Extirpated state"
#
# In the canonical OOXML this shows up as the leading "In 20% " text
# being split into its own run (leaving the remainder of the sentence
# as a second, pre-existing run) so that just the "20%" -> "21%" prefix
# changes. Re-create that by editing only the first 7 characters
# ("In 20% ") of the text range in place; PowerPoint's object model
# splits the paragraph into a new run for the edited span automatically.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(12)
$shape = $slide.Shapes.Item("Rectangle 4")

$textRange = $shape.TextFrame.TextRange
$prefix = $textRange.Characters(1, 7)
if ($prefix.Text -eq "In 20% ") {
    $prefix.Text = "In 21% "
}
